$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") holds a date serial number (45182 = 2023-09-13).
# Update every data row (2 through 264) to the new date serial 45184 (2023-09-15).
$ws.Range("C2:C264").Value = 45184
